$wb = $excel.ActiveWorkbook

# --- Sheet1 (TOP): append the new related-queries data blocks (rows 27-74) ---
$ws1 = $wb.Worksheets.Item("TOP")

$data = New-Object 'object[,]' 48,2
$data[0,0] = "ginasio"
$data[0,1] = 100
$data[1,0] = "element"
$data[1,1] = 32
$data[2,0] = "element ginásio"
$data[2,1] = 31
$data[3,0] = "ginásio lisboa"
$data[3,1] = 27
$data[4,0] = "ginásio fitness hut"
$data[4,1] = 25
$data[5,0] = "ginásio perto de mim"
$data[5,1] = 22
$data[6,0] = "ginásios"
$data[6,1] = 17
$data[7,0] = "supera"
$data[7,1] = 10
$data[8,0] = "ginásio venda nova"
$data[8,1] = 7
$data[9,0] = "rpm"
$data[9,1] = 4
$data[10,0] = "ginásio element odivelas"
$data[10,1] = 3
$data[11,0] = "ginasio supera"
$data[11,1] = 3
$data[12,0] = "ginásio clube portugues"
$data[12,1] = 2
$data[13,0] = "ginasios coimbra"
$data[13,1] = 2
$data[14,0] = "ginasio clube portugues"
$data[14,1] = 1
$data[15,0] = "body balance"
$data[15,1] = 1
$data[16,0] = "ginasio"
$data[16,1] = 100
$data[17,0] = "solinca"
$data[17,1] = 57
$data[18,0] = "ginasios lisboa"
$data[18,1] = 53
$data[19,0] = "fitness hut"
$data[19,1] = 45
$data[20,0] = "ginasios porto"
$data[20,1] = 34
$data[21,0] = "fitness up"
$data[21,1] = 31
$data[22,0] = "element"
$data[22,1] = 27
$data[23,0] = "holmes place"
$data[23,1] = 24
$data[24,0] = "ginasios coimbra"
$data[24,1] = 20
$data[25,0] = "ginasios gaia"
$data[25,1] = 19
$data[26,0] = "ginasios braga"
$data[26,1] = 17
$data[27,0] = "ginasios perto de mim"
$data[27,1] = 16
$data[28,0] = "jp ginasios"
$data[28,1] = 16
$data[29,0] = "ginasios portugal"
$data[29,1] = 15
$data[30,0] = "ginasios aveiro"
$data[30,1] = 14
$data[31,0] = "ginasios leiria"
$data[31,1] = 13
$data[32,0] = "ginasios almada"
$data[32,1] = 11
$data[33,0] = "ginasios setubal"
$data[33,1] = 11
$data[34,0] = "ginasios funchal"
$data[34,1] = 10
$data[35,0] = "solinca preços"
$data[35,1] = 10
$data[36,0] = "ginasios odivelas"
$data[36,1] = 9
$data[37,0] = "ginasios em lisboa"
$data[37,1] = 8
$data[38,0] = "ginasios oeiras"
$data[38,1] = 8
$data[39,0] = "holmes place preços"
$data[39,1] = 5
$data[40,0] = "fitness hut leiria"
$data[40,1] = 4
$data[41,0] = "ginasio"
$data[41,1] = 100
$data[42,0] = "ginásios lisboa"
$data[42,1] = 91
$data[43,0] = "ginasios"
$data[43,1] = 75
$data[44,0] = "ginásios porto"
$data[44,1] = 54
$data[45,0] = "kalorias"
$data[45,1] = 20
$data[46,0] = "ginásios perto de mim"
$data[46,1] = 18
$data[47,0] = "holmes place preços"
$data[47,1] = 11

$ws1.Range("A27:B74").Value = $data

# Update the view state on the TOP sheet: it becomes the active/selected tab,
# scrolled down, with a new selection.
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 49
$ws1.Range("F68").Select()

# --- Sheet2 (RISING): no longer the selected tab ---
$ws2 = $wb.Worksheets.Item("RISING")
$ws2.Range("G18").Select()

# Re-activate TOP sheet last so it ends up as the active tab.
$ws1.Activate()
